# tarefas.docx - apply commit "Instalacao do Redux; Criacao do fetch da lista e impressao"
#
# Content-level changes (the diff's many proofErr/run-split hunks are just
# Word's live spell-checker re-tokenizing runs around individual words -
# no visible text changes result from those, so they are not reproduced
# here): only three things actually change the document:
#   1. "*Enviar objeto do produto para o localStorage" gains the same blue
#      (00B0F0) font color already used by the two paragraphs after it.
#   2. Three new task paragraphs are inserted right after the
#      "*checar unexpected use of comma operator no reduce" paragraph and
#      right before "*Criar mensagens de erro nos inputs do form":
#        *Inserir possibilidade de excluir produto
#        *Inserir metodo para controlar a quantidade
#        *Inserir metodo para ordenar lista
#   3. Nothing else in the document's visible text or paragraph formatting
#      changes.

$d = $word.ActiveDocument

# --- 1. Color "*Enviar objeto do produto para o localStorage" blue (00B0F0) ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("*Enviar objeto do produto para o localStorage")) {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Font.Color = 15773696   # wdColor BGR for RGB 00B0F0
}

# --- 2. Insert the three new paragraphs after the "*checar..." paragraph ---
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("*checar unexpected use of comma operator no reduce")) {
        $anchor = $p
        break
    }
}

$newTasks = @(
    "*Inserir possibilidade de excluir produto",
    "*Inserir método para controlar a quantidade",
    "*Inserir método para ordenar lista "
)

if ($anchor -ne $null) {
    $previous = $anchor
    foreach ($taskText in $newTasks) {
        $previous.Range.InsertParagraphAfter()
        $newIndex = $previous.Index + 1
        $p = $d.Paragraphs.Item($newIndex)
        $p.Range.Text = $taskText
        $p.Range.LanguageID = "pt-BR"
        $previous = $p
    }
}
